{"js": "// \"Agregando gestion de control\"\n// Adds a new bullet item \"Gesti\u00f3n de control.\" right before the existing\n// \"Manejo de Git y Github\" bullet (Otros conocimientos list).\n\n// Word keeps a single \"_GoBack\" bookmark around the last edited spot; the\n// author's save cleared it, so drop it too (harmless if absent).\ntry {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n} catch (e) {\n  // no _GoBack bookmark present - nothing to clean up\n}\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet gitParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Manejo de Git y Github\") {\n    gitParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (gitParagraph) {\n  // New bullet, inheriting the same list/paragraph formatting as the\n  // \"Manejo de Git y Github\" item it sits above.\n  gitParagraph.insertParagraph(\"Gesti\u00f3n de control.\", Word.InsertLocation.before);\n  await context.sync();\n\n  // Tidy up the \"Manejo de Git y Github\" run so it reads as a single run\n  // (it previously was split around a bookmark/spell-check markers).\n  const wholeRange = gitParagraph.getRange(\"Whole\");\n  const matches = wholeRange.search(\"Manejo de Git\", { matchCase: true });\n  matches.load(\"items\");\n  await context.sync();\n\n  if (matches.items.length > 0) {\n    const firstRun = matches.items[0];\n    firstRun.insertText(\"Manejo de Git y Github\", Word.InsertLocation.replace);\n    await context.sync();\n\n    const afterFirstRun = firstRun.getRange(\"After\");\n    const paragraphEnd = gitParagraph.getRange(\"End\");\n    const remainder = afterFirstRun.expandTo(paragraphEnd);\n    remainder.delete();\n    await context.sync();\n  }\n}\n", "ps1": "# \"Agregando gestion de control\"\n# Adds a new bullet item \"Gestion de control.\" right before the existing\n# \"Manejo de Git y Github\" bullet (Otros conocimientos list).\n\n$d = $word.ActiveDocument\n\n# Word keeps a single \"_GoBack\" bookmark around the last edited spot; the\n# author's save cleared it, so drop it too (harmless if absent).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# Locate the paragraph that currently reads \"Manejo de Git y Github\".\n$gitParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($text -eq \"Manejo de Git y Github\") {\n        $gitParagraph = $p\n        break\n    }\n}\n\nif ($gitParagraph -ne $null) {\n    # New bullet, inheriting the same list/paragraph formatting as the\n    # \"Manejo de Git y Github\" item it sits above. Note: once the new empty\n    # paragraph is inserted before it, $gitParagraph itself now refers to\n    # that new (empty) paragraph, and the original text moves to .Next().\n    $gitParagraph.Range.InsertParagraphBefore()\n    $gitParagraph.Range.Text = \"Gesti\u00f3n de control.\"\n    $gitParagraph = $gitParagraph.Next()\n\n    # Tidy up the \"Manejo de Git y Github\" text so it reads as a single run\n    # (it previously was split around a bookmark/spell-check markers).\n    $gitRange = $gitParagraph.Range\n    $find = $gitRange.Find\n    $find.ClearFormatting()\n    $find.Text = \"Manejo de Git y Github\"\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = \"Manejo de Git y Github\"\n    $find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2) | Out-Null\n}\n"}
